$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Total Paid : 5.34" -> "Total Paid : 5.16"
$ws.Range("D10").Value = "Total Paid : 5.16"
$ws.Range("F10").Value = 5.16

# Row 13: "Total Paid : 5.34" -> "Total Paid : 5.16"
$ws.Range("D13").Value = "Total Paid : 5.16"
$ws.Range("H13").Value = 5.16

# Row 17: new cell F17 = "Invoice Split-1"
$ws.Range("F17").Value = "Invoice Split-1"

# Row 20: "Balance : 0 of 10.68" -> "Balance : 0 of 12.14"; "Total Paid : 10.68" -> "Total Paid : 12.14"
# (set before D18/D25 so new shared-string indices are allocated in the expected order)
$ws.Range("D20").Value = "Balance : 0 of 12.14"
$ws.Range("E20").Value = "Total Paid : 12.14"

# Row 18: "Invoice Split-1" -> "Invoice Split-2"
$ws.Range("D18").Value = "Invoice Split-2"

# Row 25: "Invoice Split-1" -> "Invoice Split-2"
$ws.Range("D25").Value = "Invoice Split-2"

# Sheet view changes: topLeftCell A25 -> A8, selection B27 -> E27
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("E27").Select()
